$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 284, shifting existing rows 284-417 down to 285-418.
$ws.Rows.Item(284).Insert()

# Populate the newly inserted row 284 with the new record's data.
$ws.Range("A284").Value = 3
$ws.Range("B284").Value = "Femacal de La Calera"
$ws.Range("C284").Value = "Coquimbo"
$ws.Range("D284").Value = 44726
$ws.Range("E284").Value = 5
$ws.Range("F284").Value = "Fruta"
$ws.Range("G284").Value = 100108
$ws.Range("H284").Value = "Tropicales y subtropicales"
$ws.Range("I284").Value = 100108002
$ws.Range("J284").Value = "Mango"
$ws.Range("K284").Value = "Sin especificar"
$ws.Range("L284").Value = "Primera"
$ws.Range("M284").Value = 228
$ws.Range("N284").Value = 10000
$ws.Range("O284").Value = 10000
$ws.Range("P284").Value = 10000
$ws.Range("Q284").Value = "`$/bandeja 4 kilos"
$ws.Range("R284").Value = "Brasil"
$ws.Range("S284").Value = 2500
$ws.Range("T284").Value = 4

# Make sure the D284 cell keeps the same date-style formatting used by the
# rest of the column (style index 2, format "YYYY-MM-DD HH:MM:SS").
$ws.Range("D284").NumberFormat = $ws.Range("D285").NumberFormat
